# Timesheet update: add the "10-24-15 to 11-04-15" sheet with the next
# two weeks of entries, and make it the active tab (matching the author's
# commit: new PCB-footprint / connector-placement work logged).

$wb = $excel.ActiveWorkbook

# Keep a stable handle on the original (first) sheet, then insert the new
# sheet immediately After it, so sheet positions/handles don't shuffle.
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Add([System.Type]::Missing, $ws1)
$ws2.Name = "10-24-15 to 11-04-15"

# --- Column widths (match sheet 1's "Date" / "Description of Work" columns) ---
$ws2.Columns.Item(1).ColumnWidth = 9.8
$ws2.Columns.Item(2).ColumnWidth = 61

# --- Header row ---
$ws2.Range("A1").Value = "Date "
$ws2.Range("B1").Value = "Description of Work"
$ws2.Range("C1").Value = "Hours"

# --- Data rows ---
# row, date-serial, description, hours
$rows = @(
    @(2, 42301, "Started Footprint creation for P1 and other IC's", 1),
    @(3, 42303, "Pad design for most footprints", 1.5),
    @(4, 42307, "Initial Footprint creation completed/ Design Rules Check", 1),
    @(5, 42309, "Connector Placement on PCB. Set up design rules for placement", 0.5),
    @(6, 42311, "Fixed outstanding schematic issues from intial review", 0.75),
    @(7, 42312, "Audio, Power, Charger Circuitry placed on PCB. Updated connector placement", 1.5)
)

# The author typed row 7's description before row 6's (it was filled in out
# of order), so the shared-string table gained "Audio, Power, ..." ahead of
# "Fixed outstanding ...". Write the B-column text in that same order first
# so new shared-string entries land at matching indices, then fill in
# everything else in normal row order.
$ws2.Range("B2").Value = "Started Footprint creation for P1 and other IC's"
$ws2.Range("B3").Value = "Pad design for most footprints"
$ws2.Range("B4").Value = "Initial Footprint creation completed/ Design Rules Check"
$ws2.Range("B5").Value = "Connector Placement on PCB. Set up design rules for placement"
$ws2.Range("B7").Value = "Audio, Power, Charger Circuitry placed on PCB. Updated connector placement"
$ws2.Range("B6").Value = "Fixed outstanding schematic issues from intial review"

foreach ($r in $rows) {
    $rowIdx = $r[0]
    $ws2.Range("A$rowIdx").Value = $r[1]
    $ws2.Range("A$rowIdx").NumberFormat = "d-mmm"
    $ws2.Range("A$rowIdx").VerticalAlignment = -4108   # xlVAlignCenter

    $ws2.Range("B$rowIdx").WrapText = $true

    $ws2.Range("C$rowIdx").Value = $r[3]
}

# Last entry ("Audio, Power, Charger...") wraps to two lines like the
# equivalent long rows on sheet 1.
$ws2.Rows.Item(7).RowHeight = 30

# --- Total row ---
$ws2.Range("A16").Value = "Total:"
$ws2.Range("A16:B16").Merge()
$ws2.Range("A16:B16").Font.Bold = $true
$ws2.Range("A16:B16").HorizontalAlignment = -4108   # xlCenter
$ws2.Range("C16").Formula = "=SUM(C2:C7)"

# --- Page margins (match sheet 1) ---
$ws2.PageSetup.LeftMargin = $excel.InchesToPoints(0.7)
$ws2.PageSetup.RightMargin = $excel.InchesToPoints(0.7)
$ws2.PageSetup.TopMargin = $excel.InchesToPoints(0.75)
$ws2.PageSetup.BottomMargin = $excel.InchesToPoints(0.75)
$ws2.PageSetup.HeaderMargin = $excel.InchesToPoints(0.3)
$ws2.PageSetup.FooterMargin = $excel.InchesToPoints(0.3)

# --- Selection / activation: new sheet becomes the active tab, cursor on B6 ---
$ws2.Range("B6").Select()
$ws2.Activate()
